# Scheduled-runner style refresh of cached market-board figures in the
# per-job "Leve Profits" tables (columns H-N: currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)). Mirrors an automated data pull that
# overwrites previously cached values with freshly retrieved ones.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 608.875
$ws.Range("I18").Value = 608.875
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 608.875
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -324.875
$ws.Range("N18").ClearContents()

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 148.95
$ws.Range("I28").Value = 155.5
$ws.Range("J28").Value = 139.125
$ws.Range("K28").Value = 155.5
$ws.Range("L28").Value = 139.125
$ws.Range("M28").Value = 329.5
$ws.Range("N28").Value = -1109.125

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 10418378
$ws.Range("I100").Value = 10418378
$ws.Range("K100").Value = 10418378
$ws.Range("M100").Value = -10417837

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 11364363
$ws.Range("I107").Value = 20833724
$ws.Range("J107").Value = 1129.5
$ws.Range("K107").Value = 20833724
$ws.Range("L107").Value = 1129.5
$ws.Range("M107").Value = -20831804
$ws.Range("N107").Value = -4969.5

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 92815.45
$ws.Range("I111").Value = 1981.75
$ws.Range("J111").Value = 144720.42
$ws.Range("K111").Value = 5945.25
$ws.Range("L111").Value = 434161.26
$ws.Range("M111").Value = -2878.25
$ws.Range("N111").Value = -440295.26

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 5904
$ws.Range("I116").Value = 7211.4707
$ws.Range("J116").Value = 3125.625
$ws.Range("K116").Value = 7211.4707
$ws.Range("L116").Value = 3125.625
$ws.Range("M116").Value = -3769.4707
$ws.Range("N116").Value = -10009.625

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 7938010.5
$ws.Range("I132").Value = 940.5862
$ws.Range("J132").Value = 25643782
$ws.Range("K132").Value = 2821.7586
$ws.Range("L132").Value = 76931346
$ws.Range("M132").Value = -291.7586000000001
$ws.Range("N132").Value = -76936406

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 1794.434
$ws.Range("I135").Value = 1656.2778
$ws.Range("J135").Value = 2087
$ws.Range("K135").Value = 14906.5002
$ws.Range("L135").Value = 18783
$ws.Range("M135").Value = -12371.5002
$ws.Range("N135").Value = -23853

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1481.5526
$ws.Range("I137").Value = 1135.3
$ws.Range("J137").Value = 2780
$ws.Range("K137").Value = 3405.9
$ws.Range("L137").Value = 8340
$ws.Range("M137").Value = -855.8999999999996
$ws.Range("N137").Value = -13440

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4239.1167
$ws.Range("I138").Value = 806.87805
$ws.Range("J138").Value = 11645.526
$ws.Range("K138").Value = 2420.63415
$ws.Range("L138").Value = 34936.578
$ws.Range("M138").Value = 2719.36585
$ws.Range("N138").Value = -45216.578

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 1110.591
$ws.Range("I141").Value = 1110.591
$ws.Range("K141").Value = 3331.773
$ws.Range("M141").Value = 1848.227

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 9616562
$ws.Range("I74").Value = 1056.5278
$ws.Range("J74").Value = 31251448
$ws.Range("K74").Value = 1056.5278
$ws.Range("L74").Value = 31251448
$ws.Range("M74").Value = -182.5278000000001
$ws.Range("N74").Value = -31253196

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 9616562
$ws.Range("I77").Value = 1056.5278
$ws.Range("J77").Value = 31251448
$ws.Range("K77").Value = 5282.639
$ws.Range("L77").Value = 156257240
$ws.Range("M77").Value = -914.6390000000001
$ws.Range("N77").Value = -156265976

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2180720.8
$ws.Range("I102").Value = 2471216.8
$ws.Range("K102").Value = 2471216.8
$ws.Range("M102").Value = -2469594.8

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2922.2646
$ws.Range("I132").Value = 1789.6666
$ws.Range("J132").Value = 5640.5
$ws.Range("K132").Value = 5368.9998
$ws.Range("L132").Value = 16921.5
$ws.Range("M132").Value = -2838.9998
$ws.Range("N132").Value = -21981.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 2444.4443

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 66667800
$ws.Range("I99").Value = 90910136
$ws.Range("K99").Value = 90910136
$ws.Range("M99").Value = -90908638

# Row 103 (Leve Item ID 18514)
$ws.Range("H103").Value = 36500
$ws.Range("J103").Value = 36500
$ws.Range("L103").Value = 36500
$ws.Range("N103").Value = -38844

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 33335316
$ws.Range("I105").Value = 50001616
$ws.Range("K105").Value = 50001616
$ws.Range("M105").Value = -49999869

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3449961
$ws.Range("I132").Value = 4167877.8
$ws.Range("J132").Value = 3960.2
$ws.Range("K132").Value = 12503633.4
$ws.Range("L132").Value = 11880.6
$ws.Range("M132").Value = -12501103.4
$ws.Range("N132").Value = -16940.6

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 4976806.5
$ws.Range("I134").Value = 6291039.5
$ws.Range("K134").Value = 18873118.5
$ws.Range("M134").Value = -18870583.5

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1306.2858
$ws.Range("I68").Value = 499
$ws.Range("K68").Value = 1497
$ws.Range("M68").Value = -686

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1306.2858
$ws.Range("I71").Value = 499
$ws.Range("K71").Value = 4491
$ws.Range("M71").Value = -435

# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 12500717
$ws.Range("I98").Value = 103
$ws.Range("J98").Value = 14286519
$ws.Range("K98").Value = 309
$ws.Range("L98").Value = 42859557
$ws.Range("M98").Value = 1189
$ws.Range("N98").Value = -42862553

# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 417242.16
$ws.Range("I113").Value = 629.5714
$ws.Range("J113").Value = 1000499.8
$ws.Range("K113").Value = 1888.7142
$ws.Range("L113").Value = 3001499.4
$ws.Range("M113").Value = 281.2857999999999
$ws.Range("N113").Value = -3005839.4

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 5046.4585
$ws.Range("I122").Value = 475
$ws.Range("K122").Value = 4275
$ws.Range("M122").Value = -1825

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 27779930
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 27779930
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 250019370
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -250024430

$ws = $wb.Worksheets.Item("GSM")
# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 9404.571
$ws.Range("I57").Value = 4177.5
$ws.Range("J57").Value = 10275.75
$ws.Range("K57").Value = 4177.5
$ws.Range("L57").Value = 10275.75
$ws.Range("M57").Value = -3357.5
$ws.Range("N57").Value = -11915.75

# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 628.8889
$ws.Range("I97").Value = 620
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 620
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -124
$ws.Range("N97").Value = -1692

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 43479588
$ws.Range("I113").Value = 76924170
$ws.Range("K113").Value = 76924170
$ws.Range("M113").Value = -76922000

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 9756658
$ws.Range("I132").Value = 12392528
$ws.Range("J132").Value = 3939.3
$ws.Range("K132").Value = 37177584
$ws.Range("L132").Value = 11817.9
$ws.Range("M132").Value = -37175054
$ws.Range("N132").Value = -16877.9

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 8069.684
$ws.Range("I136").Value = 6105.9614
$ws.Range("J136").Value = 12324.417
$ws.Range("K136").Value = 18317.8842
$ws.Range("L136").Value = 36973.251
$ws.Range("M136").Value = -15767.8842
$ws.Range("N136").Value = -42073.251

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1170.9718
$ws.Range("I132").Value = 818.10205
$ws.Range("K132").Value = 2454.30615
$ws.Range("M132").Value = 75.69385000000011

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 3404354.5
$ws.Range("I136").Value = 4000.9614
$ws.Range("J136").Value = 7248232
$ws.Range("K136").Value = 12002.8842
$ws.Range("L136").Value = 21744696
$ws.Range("M136").Value = -9452.8842
$ws.Range("N136").Value = -21749796
